# Auto-generated edit script: updates the cryptos price/volume table
# to the refreshed snapshot values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = '60.893.67'
$ws.Range("E2").Value = '  -0.69%  '

# --- Row 3 ---
$ws.Range("D3").Value = '3.379.22'
$ws.Range("E3").Value = '  +1.06%  '

# --- Row 4 ---
$ws.Range("E4").Value = '  -0.09%  '

# --- Row 5 ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '403.85'
$ws.Range("E5").Value = '  -2.04%  '

# --- Row 6 ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.62'
$ws.Range("E6").Value = '  +14.00%  '

# --- Row 7 ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.596'
$ws.Range("E7").Value = '  +4.82%  '

# --- Row 8 ---
$ws.Range("D8").Value = '3.371.52'
$ws.Range("E8").Value = '  +1.03%  '

# --- Row 9 ---
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.08%  '

# --- Row 10 ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.668'
$ws.Range("E10").Value = '  +6.77%  '

# --- Row 11 ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +15.17%  '

# --- Row 12 ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.03'
$ws.Range("E12").Value = '  +7.91%  '

# --- Row 13 ---
$ws.Range("E13").Value = '  -1.05%  '

# --- Row 14 ---
$ws.Range("D14").Value = '3.915.15'
$ws.Range("E14").Value = '  +0.46%  '

# --- Row 15 ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.51'
$ws.Range("E15").Value = '  +3.13%  '

# --- Row 16 ---
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.56'
$ws.Range("E16").Value = '  +2.24%  '

# --- Row 17 ---
$ws.Range("D17").Value = '3.397.30'
$ws.Range("E17").Value = '  -0.68%  '

# --- Row 18 ---
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.54'
$ws.Range("E18").Value = '  +9.53%  '

# --- Row 19 ---
$ws.Range("D19").Value = '60.887.61'
$ws.Range("E19").Value = '  -0.33%  '

# --- Row 20 ---
$ws.Range("E20").Value = '  +0.26%  '

# --- Row 21 ---
$ws.Range("E21").Value = '  +15.70%  '

# --- Row 22 ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.22'
$ws.Range("E22").Value = '  -1.19%  '

# --- Row 23 ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.10'
$ws.Range("E23").Value = '  +10.90%  '

# --- Row 24 ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.93'
$ws.Range("E24").Value = '  +4.78%  '

# --- Row 25 ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '305.62'
$ws.Range("E25").Value = '  +2.48%  '

# --- Row 26 ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.14'
$ws.Range("E26").Value = '  +1.26%  '

# --- Row 27 ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.60'
$ws.Range("E27").Value = '  +13.92%  '

# --- Row 28 ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '29.44'
$ws.Range("E28").Value = '  +2.36%  '

# --- Row 29 ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.62'
$ws.Range("E29").Value = '  +3.00%  '

# --- Row 30 ---
$ws.Range("E30").Value = '  -0.22%  '

# --- Row 31 ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.172'
$ws.Range("E31").Value = '  +1.54%  '

# --- Row 32 ---
$ws.Range("E32").Value = '  +3.81%  '

# --- Row 33 ---
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.62'
$ws.Range("E33").Value = '  +3.11%  '

# --- Row 34 ---
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.72'
$ws.Range("E34").Value = '  +8.06%  '

# --- Row 35 ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.56'
$ws.Range("E35").Value = '  +5.40%  '

# --- Row 36 ---
$ws.Range("E36").Value = '  +0.16%  '

# --- Row 37 ---
$ws.Range("E37").Value = '  +1.09%  '

# --- Row 38 ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.08'
$ws.Range("E38").Value = '  +0.18%  '

# --- Row 39 ---
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.997'
$ws.Range("E39").Value = '  -0.33%  '

# --- Row 40 ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  +2.19%  '

# --- Row 41 ---
$ws.Range("E41").Value = '  -4.54%  '

# --- Row 42 ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.124'
$ws.Range("E42").Value = '  +3.40%  '

# --- Row 43 ---
$ws.Range("E43").Value = '  +4.03%  '

# --- Row 44 ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '134.70'
$ws.Range("E44").Value = '  -3.13%  '

# --- Row 45 ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.282'
$ws.Range("E45").Value = '  +1.05%  '

# --- Row 46 ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.86'
$ws.Range("E46").Value = '  +3.09%  '

# --- Row 47 ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.89'
$ws.Range("E47").Value = '  +1.61%  '

# --- Row 48 ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.22'
$ws.Range("E48").Value = '  +0.91%  '

# --- Row 49 ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.78'
$ws.Range("E49").Value = '  +2.80%  '

# --- Row 50 ---
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '3.726.58'
$ws.Range("E50").Value = '  -3.75%  '

# --- Row 51 ---
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.134.23'
$ws.Range("E51").Value = '  -0.42%  '
